# "Generate Report for Handoff" - regenerate the localization-status report:
# the zh-cn / de-de items move from "In Translation" to "Ready for handoff"
# and the handoff timestamps are refreshed. Column widths are re-fitted to
# the (now longer) status text, same as the reporting tool would do when it
# re-writes the workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # zh-cn detail sheet status
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # de-de detail sheet status

# --- Refresh the handoff timestamps -----------------------------------------
$wsOverview.Range("G2").Value = "2016-08-27 10:39:34" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-08-27 10:39:34" # de-de Latest Handoff Datetime
$wsZhCn.Range("H2").Value     = "2016-08-27 10:39:30" # zh-cn Latest Handoff Datetime

# --- Re-fit the Status columns now that the text is longer ------------------
$wsOverview.Columns("E:E").ColumnWidth = 16.29
$wsOverview.Columns("F:F").ColumnWidth = 16.29
$wsZhCn.Columns("C:C").ColumnWidth = 16.29
$wsDeDe.Columns("C:C").ColumnWidth = 16.29
